$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell F1 with the same style as E1 (bold/centered/bordered header style)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"
$excel.CutCopyMode = $false

# Populate time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 13:38:47.753561"
$ws.Range("F3").Value = "2021-10-05 13:38:47.753571"
$ws.Range("F4").Value = "2021-10-05 13:38:47.753573"
$ws.Range("F5").Value = "2021-10-05 13:38:47.753575"
$ws.Range("F6").Value = "2021-10-05 13:38:47.753578"
$ws.Range("F7").Value = "2021-10-05 13:38:47.753580"
$ws.Range("F8").Value = "2021-10-05 13:38:47.753582"
$ws.Range("F9").Value = "2021-10-05 13:38:47.753584"
$ws.Range("F10").Value = "2021-10-05 13:38:47.753586"
$ws.Range("F11").Value = "2021-10-05 13:38:47.753588"
$ws.Range("F12").Value = "2021-10-05 13:38:47.753590"
$ws.Range("F13").Value = "2021-10-05 13:38:47.753592"
$ws.Range("F14").Value = "2021-10-05 13:38:47.753594"
$ws.Range("F15").Value = "2021-10-05 13:38:47.753596"
$ws.Range("F16").Value = "2021-10-05 13:38:47.753598"
$ws.Range("F17").Value = "2021-10-05 13:38:47.753600"
$ws.Range("F18").Value = "2021-10-05 13:38:47.753603"
$ws.Range("F19").Value = "2021-10-05 13:38:47.753605"
$ws.Range("F20").Value = "2021-10-05 13:38:47.753607"
$ws.Range("F21").Value = "2021-10-05 13:38:47.753609"
$ws.Range("F22").Value = "2021-10-05 13:38:47.753610"
$ws.Range("F23").Value = "2021-10-05 13:38:47.753612"
$ws.Range("F24").Value = "2021-10-05 13:38:47.753614"
$ws.Range("F25").Value = "2021-10-05 13:38:47.753616"
$ws.Range("F26").Value = "2021-10-05 13:38:47.753618"
$ws.Range("F27").Value = "2021-10-05 13:38:47.753620"
$ws.Range("F28").Value = "2021-10-05 13:38:47.753622"
$ws.Range("F29").Value = "2021-10-05 13:38:47.753624"
$ws.Range("F30").Value = "2021-10-05 13:38:47.753626"
